# Ajuste envio de Correos
# - Renumber the "consecutivo" column (A) for the three data rows.
# - Unify the "nombre_producto" column (B) to "ortografia" for all rows.
# - Fix a typo in Janluy's e-mail address (underscore -> dot) while keeping
#   the existing mailto: hyperlink target untouched.
# - Widen columns F and O to fit their (new) contents.
# - Leave the active selection on the (now interesting) e-mail cell O3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Janluy Leonard Moreno Coronado) ---------------------------
$ws.Range("A2").Value = "00001"
$ws.Range("B2").Value = "ortografia"
$ws.Range("O2").Value = "janluy.moreno@cun.edu.co"

# --- Row 3 (Nelson Andres Cardenas Velasquez) --------------------------
$ws.Range("A3").Value = "00002"
$ws.Range("B3").Value = "ortografia"

# --- Row 4 (johan camilo triana avendaño) -------------------------------
$ws.Range("A4").Value = "00003"
# B4 previously had no explicit cell style; after the edit it matches the
# other rows' "nombre_producto" cells (style index 2), so copy that
# formatting across before overwriting the text.
$ws.Range("C2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "ortografia"

# --- Column widths (best-fit for the new long values) -------------------
$ws.Columns("F").ColumnWidth = 13
$ws.Columns("O").ColumnWidth = 27

# --- Active selection ---------------------------------------------------
$ws.Range("O3").Select()
